# Site updated: 2020-04-24 15:29:01
#
# Adds two new file-note entries (rows 72 and 73) to Sheet1, and restores
# the view state (selection / scroll position) that results from having
# scrolled down to see the newly appended rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data rows -------------------------------------------------------
$ws.Range("A72").Value = "UML-SequenceDiagram"
$ws.Range("B72").Value = "UML学习笔记——顺序图（时序图）"

$ws.Range("A73").Value = "find-similar-string"
$ws.Range("B73").Value = "找相似串问题——算法设计课程"

# --- View / window state --------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 54
$win.ScrollColumn = 1

$ws.Range("B60").Select()
